$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Preserve the "yellow highlight" cell format before values/formats
# in column B (rows 2,4,5,6,7) get overwritten, by stashing a copy
# of it onto an unused scratch cell first.
# -----------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Apply the stashed yellow highlight format onto the new "label"
# column (A) for the dataset-level metadata rows.
$ws.Range("F1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)

# Reset the old label column (B) for those rows back to the plain
# (non-highlighted) format, using an already-plain cell as source.
$ws.Range("A3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)

# Remove the scratch cell contents/format.
$ws.Range("F1").Clear()

# -----------------------------------------------------------------
# Rearrange table contents.
# Old layout: A=Identifier(Dataset/Column)  B=meta_data/value label  C=Value(number)/Detail(json)
# New layout: A=meta_data/value label       B=Value(number/short)    C=Identifier(dataset/column)  D=Detail(json)
# -----------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "meta_data"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Identifier"
$ws.Range("D1").Value = "Detail"

# Dataset-level metadata rows (2-7)
$ws.Range("A2").Value = "# of columns"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "dataset"

$ws.Range("A3").Value = "# of rows"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = "dataset"

$ws.Range("A4").Value = "# of Text columns"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "dataset"

$ws.Range("A5").Value = "# of Numeric Columns"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "dataset"

$ws.Range("A6").Value = "# of Boolean Columns"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "dataset"

$ws.Range("A7").Value = "# of Date Columns"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "dataset"

# Column-level metadata rows (8-10)
$ws.Range("A8").Value = "Text 1"
$ws.Range("B8").Value = "STRING"
$ws.Range("C8").Value = "column"
$ws.Range("D8").Value = "{`n    column_name: 'Text 1',`n    column_type: 'STRING',`n    count: 45,`n    unique: 39,`n    top: 'a',`n    value_counts_top_10: {`n        'a': 20,`n        'b': 40,`n    },`n    null_count: 45`n}"

$ws.Range("A9").Value = "Numeric 1"
$ws.Range("B9").Value = "INT64"
$ws.Range("C9").Value = "column"
$ws.Range("D9").Value = "{`n    column_name: 'Numeric 1',`n    column_type: 'INT64',`n    count: 45,`n    mean: 39,`n    std: 39,`n    min: 39,`n    max: 39,`n    25%: 39,`n    50%: 39,`n    75%: 39,`n    null_count: 45`n}"

$ws.Range("A10").Value = "Date 1"
$ws.Range("B10").Value = "DATETIME"
$ws.Range("C10").Value = "column"
$ws.Range("D10").Value = "{`n    column_name: 'Date 1',`n    column_type: 'datetime',`n    count: 45,`n    min: 39,`n    max: 39,`n    estimated_granularity: 'day',`n    null_count: 45,`n    missing_dates: []`n}"

# -----------------------------------------------------------------
# Column widths: swap widths of columns 1 and 2, column 3 reverts
# to the workbook default width.
# -----------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.7109375
$ws.Columns.Item(2).ColumnWidth = 18.140625
$ws.Columns.Item(3).ColumnWidth = 9.140625

# -----------------------------------------------------------------
# Update the active selection.
# -----------------------------------------------------------------
$ws.Range("B8").Select()
